$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# D-column values are forced to remain plain text (matching the source feed
# formatting, e.g. "644.10" or "69.439.34") by briefly applying a Text number
# format before assignment, then restoring the default "Normal" style so the
# cell keeps its original (unstyled) appearance.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.439.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.13%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.677.52"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.27%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "644.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.48%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.498"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.73%  "

$ws.Range("E9").Value = "  -0.41%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.78%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.447"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.15%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000232"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.40%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.296.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.85%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.671.25"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.42%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.418.57"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "16.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.13%  "

$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "465.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.36%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.56%  "

$ws.Range("E22").Value = "  -0.92%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "79.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.62%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.823.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.27%  "

$ws.Range("E25").Value = "  -0.03%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000126"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.12%  "

$ws.Range("E27").Value = "  +0.19%  "

$ws.Range("E28").Value = "  -0.68%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.63"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.32%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.72"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.86%  "

$ws.Range("E31").Value = "  +0.81%  "

$ws.Range("E32").Value = "  -0.16%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.17%  "

$ws.Range("E34").Value = "  +3.48%  "

$ws.Range("E35").Value = "  -1.71%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.670.48"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.18%  "

$ws.Range("E37").Value = "  +1.50%  "

$ws.Range("E38").Value = "  -0.03%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.22%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "178.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.76%  "

$ws.Range("E41").Value = "  -0.08%  "

$ws.Range("E42").Value = "  -0.33%  "

$ws.Range("E43").Value = "  -1.91%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.927"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.66%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "46.66"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.68%  "

$ws.Range("E46").Value = "  +2.36%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "27.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.96%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000270"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.25%  "

$ws.Range("E49").Value = "  -3.44%  "

$ws.Range("E50").Value = "  +0.74%  "

$ws.Range("E51").Value = "  -3.79%  "

